# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> name="Office Theme" (classic blue "Office" palette)
#   ppt/theme/theme2.xml  -> name="Integral"     (green "Integral" palette),
#                             linked from both presentation.xml and
#                             slideMaster1.xml (i.e. the palette actually
#                             painting every slide in the deck).
#
# The commit swaps the two themes' contents: the slide master/presentation
# theme (theme2.xml) ends up holding the "Office Theme" colour scheme, while
# theme1.xml ends up holding the "Integral" colour scheme. fontScheme /
# fmtScheme are identical between the two parts both before and after, so
# the only thing that actually changes is the 12-slot colour scheme (and
# the cosmetic theme/clrScheme "name" attributes).
#
# Reproduce this through the object model by rewriting the slide master's
# theme colour scheme (the part every slide/master in this deck renders
# with) to the "Office Theme" palette values that theme1.xml used to hold.

function HexToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Theme colours, in the fixed OOXML <a:clrScheme> slot order that
# ThemeColorScheme.Item(1..12) walks: dk1, lt1, dk2, lt2, accent1..6,
# hlink, folHlink. These are the values theme1.xml ("Office Theme") had.
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToRGB $officeThemeColors[$i - 1]
}

# Keep the design's display name in step with its new palette.
$p.Designs.Item(1).Name = "Office Theme"
